$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force Excel to store the literal text rather than auto-converting a
    # numeric-looking string ("714.61", "10.72", ...) into a number: type it
    # with a leading apostrophe (the standard "force text" entry method),
    # then strip the resulting quote-prefix style back to Normal so no
    # stray formatting is left behind on the cell.
    $ws.Range($range).Value = "'" + $text
    $ws.Range($range).Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "71.134.06"
$ws.Range("E2").Value = "  +0.66%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.836.78"
$ws.Range("E3").Value = "  +1.07%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "714.61"
$ws.Range("E5").Value = "  +1.93%  "

# Row 6 - Solana
Set-TextValue "D6" "172.80"
$ws.Range("E6").Value = "  +0.31%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.836.59"
$ws.Range("E7").Value = "  +1.10%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.04%  "

# Row 9 - XRP
Set-TextValue "D9" "0.527"

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.46%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +1.47%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -0.04%  "

# Row 13 - ShibaInu
Set-TextValue "D13" "0.0000256"
$ws.Range("E13").Value = "  +0.16%  "

# Row 14 - Avalanche
Set-TextValue "D14" "36.78"
$ws.Range("E14").Value = "  +1.84%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "4.483.46"
$ws.Range("E15").Value = "  +1.08%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.953.03"
$ws.Range("E16").Value = "  +4.49%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "71.111.28"
$ws.Range("E17").Value = "  +0.75%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  +0.38%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  +0.68%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  -1.80%  "

# Row 21 - Uniswap
Set-TextValue "D21" "10.72"
$ws.Range("E21").Value = "  -4.21%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "494.68"
$ws.Range("E22").Value = "  +3.11%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.726"
$ws.Range("E23").Value = "  +2.07%  "

# Row 24 - Litecoin
Set-TextValue "D24" "85.17"
$ws.Range("E24").Value = "  +1.48%  "

# Row 25 - PEPE
$ws.Range("E25").Value = "  +2.45%  "

# Row 26 - RenderToken
Set-TextValue "D26" "10.66"
$ws.Range("E26").Value = "  +1.91%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "12.17"
$ws.Range("E27").Value = "  -1.24%  "

# Row 28 - PancakeSwap
Set-TextValue "D28" "3.21"
$ws.Range("E28").Value = "  +3.04%  "

# Row 29 - Fetch.AI
Set-TextValue "D29" "2.10"
$ws.Range("E29").Value = "  -2.07%  "

# Row 30 - Dai
$ws.Range("E30").Value = "  -0.08%  "

# Row 31 - NEARProtocol
Set-TextValue "D31" "7.50"
$ws.Range("E31").Value = "  -0.61%  "

# Row 32 - ImmutableX
Set-TextValue "D32" "2.24"
$ws.Range("E32").Value = "  -1.86%  "

# Row 33 - EthereumClassic
Set-TextValue "D33" "29.39"
$ws.Range("E33").Value = "  -0.01%  "

# Row 34 - Kaspa
$ws.Range("E34").Value = "  -4.22%  "

# Row 35 - Aptos
$ws.Range("E35").Value = "  -0.19%  "

# Row 36 - RenzoRestakedETH
Set-TextValue "D36" "3.800.78"
$ws.Range("E36").Value = "  +1.55%  "

# Row 37 - Binance-PegBSC-USD
Set-TextValue "D37" "0.998"
$ws.Range("E37").Value = "  -0.22%  "

# Row 38 - Hedera
$ws.Range("E38").Value = "  +0.82%  "

# Row 39 - Mantle
$ws.Range("E39").Value = "  +5.88%  "

# Row 40 - Filecoin
$ws.Range("E40").Value = "  +0.56%  "

# Row 41 - dogwifhat
Set-TextValue "D41" "3.36"
$ws.Range("E41").Value = "  -0.76%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  +2.10%  "

# Row 43 - USDe
$ws.Range("E43").Value = "  +0.00%  "

# Row 44 - FirstDigitalUSD
$ws.Range("E44").Value = "  +0.16%  "

# Row 45 - FLOKI
Set-TextValue "D45" "0.000321"
$ws.Range("E45").Value = "  +0.33%  "

# Row 46 - Monero
Set-TextValue "D46" "163.47"
$ws.Range("E46").Value = "  +0.06%  "

# Row 47 - OKB
Set-TextValue "D47" "48.86"
$ws.Range("E47").Value = "  +0.03%  "

# Row 48 - Bittensor
Set-TextValue "D48" "423.48"
$ws.Range("E48").Value = "  +3.55%  "

# Row 49 - ONDO
Set-TextValue "D49" "1.39"
$ws.Range("E49").Value = "  -0.01%  "

# Row 50 - Cosmos
$ws.Range("E50").Value = "  +0.25%  "

# Row 51 - TheGraph
$ws.Range("E51").Value = "  -0.99%  "
